# Add the "債務" (debt) worksheet after "汽車", matching the layout/style
# used by the other three sheets (土地/建物/汽車), and fill it with the
# legislator's debt records.

$wb = $excel.ActiveWorkbook

# New sheet goes at the end of the tab strip (after the last existing sheet).
$ws = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws.Name = "債務"

# ---- Header row (row 1, columns B..N), left to right -----------------------
$headerCols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L", "M", "N")
$headers    = @("species", "debtor", "owner", "total", "register_date", "register_reason", "property_category", "category", "date", "legislator_name", "legislator_id", "source_file", "index")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Range($headerCols[$i] + "1").Value = $headers[$i]
}

# Header styling: bold font, thin border all around, centered (matches the
# style used for row 1 on the other sheets).
$headerRange = $ws.Range("B1:N1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

# ---- Data rows (rows 2..10) -------------------------------------------------
# Index column (A, also duplicated in N)
$indexes = @(86, 87, 88, 89, 90, 91, 92, 93, 94)
$n = $indexes.Length

# Column A ("index") gets the same bold/border/centered styling as the header
# on the other sheets.
$colARange = $ws.Range("A2:A10")
$colARange.Font.Bold = $true
$colARange.Borders.LineStyle = 1
$colARange.HorizontalAlignment = -4108
$colARange.VerticalAlignment = -4160

# Per-column values, top to bottom (rows 2..10) -- written column by column so
# new shared strings land in the same order the source workbook used.
$colB = @("房貸", "房貸", "房貸", "房貸", "房貸", "房貸", "房貸", "房貸", "借款")
$colC = @("林世嘉", "林世嘉", "蔡篤堅", "蔡篤堅", "蔡篤堅", "蔡篤堅", "蔡篤堅", "蔡篤堅", "蔡篤堅")
$colD = @("兆豐國際商業銀行台北復興分行", "兆豐國際商業銀行台北復興分行", "第一商業銀行盧州分行", "第一商業銀行盧州分行", "第一商業銀行盧州分行", "第一商業銀行盧州分行", "第一商業銀行盧州分行", "第一商業銀行盧州分行", "郭素珍")
$colE = @(5566788, 6675465, 2037093, 6739079, 2161265, 563630, 543011, 714257, 2000000)
$colF = @("97年07月", "93年02月", "97年03月", "97年04月", "98年04月", "92年01月", "92年01月", "92年01月", "92年02月")
$colG = @("購屋", "購屋", "購屋", "購屋", "購屋", "購屋", "購屋", "購屋", "資金週轉")
$colH = @("debt", "debt", "debt", "debt", "debt", "debt", "debt", "debt", "debt")
$colI = @("normal", "normal", "normal", "normal", "normal", "normal", "normal", "normal", "normal")
$colJ = @("2012-05-01", "2012-05-01", "2012-05-01", "2012-05-01", "2012-05-01", "2012-05-01", "2012-05-01", "2012-05-01", "2012-05-01")
$colK = @("林世嘉", "林世嘉", "林世嘉", "林世嘉", "林世嘉", "林世嘉", "林世嘉", "林世嘉", "林世嘉")
$colL = @(1740, 1740, 1740, 1740, 1740, 1740, 1740, 1740, 1740)
$colM = @("tmpada11", "tmpada11", "tmpada11", "tmpada11", "tmpada11", "tmpada11", "tmpada11", "tmpada11", "tmpada11")

for ($i = 0; $i -lt $n; $i++) {
    $ws.Range("A" + (2 + $i)).Value = $indexes[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Range("B" + (2 + $i)).Value = $colB[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Range("C" + (2 + $i)).Value = $colC[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Range("D" + (2 + $i)).Value = $colD[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Range("E" + (2 + $i)).Value = $colE[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Range("F" + (2 + $i)).Value = $colF[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Range("G" + (2 + $i)).Value = $colG[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Range("H" + (2 + $i)).Value = $colH[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Range("I" + (2 + $i)).Value = $colI[$i]
}

# Column J ("date") holds an ISO-looking string ("2012-05-01") which Excel
# would otherwise auto-convert to a date serial number. Force it to stay
# text by pre-formatting as Text, then strip the leftover text format so the
# cell ends up unstyled (same as its sibling data cells) once the literal
# string value has been committed.
$colJRange = $ws.Range("J2:J10")
$colJRange.NumberFormat = "@"
for ($i = 0; $i -lt $n; $i++) {
    $ws.Range("J" + (2 + $i)).Value = $colJ[$i]
}
$colJRange.ClearFormats()

for ($i = 0; $i -lt $n; $i++) {
    $ws.Range("K" + (2 + $i)).Value = $colK[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Range("L" + (2 + $i)).Value = $colL[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Range("M" + (2 + $i)).Value = $colM[$i]
}
for ($i = 0; $i -lt $n; $i++) {
    $ws.Range("N" + (2 + $i)).Value = $indexes[$i]
}

Write-Output "done"
